# Add a new column C ("t+3") next to the existing A/B columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting from B1 (bold/border/center) onto C1, then set its value.
$ws.Range("B1").Copy($ws.Range("C1"))
$ws.Range("C1").Value = 2

# Values for the new "t+3" data column.
$values = @(
    -4.993368022640259,
    -1.165564360090414,
    -0.06515441686030865,
    -0.4032515873081615,
    0.01465567179956126,
    0.1084370207011733,
    0.1341971137761105
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $values[$i]
}
